$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Shift the "Key" legend table from columns I:J to J:K (rows 2-16)
#    A new column was effectively inserted before the legend table.
# -----------------------------------------------------------------
for ($r = 16; $r -ge 2; $r--) {
    $iCell = $ws.Cells.Item($r, 9)
    $jCell = $ws.Cells.Item($r, 10)
    $iHasValue = $iCell.Value() -ne $null
    $jHasValue = $jCell.Value() -ne $null
    if ($jHasValue) {
        $ws.Cells.Item($r, 10).Copy($ws.Cells.Item($r, 11)) | Out-Null
    }
    if ($iHasValue) {
        $ws.Cells.Item($r, 9).Copy($ws.Cells.Item($r, 10)) | Out-Null
    }
}
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 9).Clear() | Out-Null
}

# -----------------------------------------------------------------
# 2. Finish the shot list: scenes 7 and 8 (rows 28-34), replacing the
#    old YouTube-reference notes that used to sit in rows 29-31.
# -----------------------------------------------------------------
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 0.002777777777777778
$ws.Range("C28").NumberFormat = "h:mm"
$ws.Range("D28").Value = "ES/LS"
$ws.Range("E28").Value = "Crane down and back"
$ws.Range("F28").Value = "Fridge"
$ws.Range("G28").Value = "fridge flies in from the distance above, bounces and creates a dust cloud, then flies off screen"

$ws.Range("A29").Value = 7
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 0.001388888888888889
$ws.Range("C29").NumberFormat = "h:mm"
$ws.Range("D29").Value = "LS"
$ws.Range("E29").Value = "Stationary (shaky cam)"
$ws.Range("F29").Value = "Fridge"
$ws.Range("G29").Value = "fridge bounces down a hill from right to left"

$ws.Range("A30").Value = 7
$ws.Range("B30").Value = 3
$ws.Range("C30").Value = 0.0020833333333333333
$ws.Range("C30").NumberFormat = "h:mm"
$ws.Range("D30").Value = "LS, HA"
$ws.Range("E30").Value = "Stationary"
$ws.Range("F30").Value = "Fridge"
$ws.Range("G30").Value = "fridge bounces down hill, away from camera"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = 4
$ws.Range("C31").Value = 0.0062499999999999995
$ws.Range("C31").NumberFormat = "h:mm"
$ws.Range("D31").Value = "LS, LA"
$ws.Range("E31").Value = "Stationary (shaky cam)"
$ws.Range("F31").Value = "Fridge, Indy"
$ws.Range("G31").Value = "fridge bounces down hill, towards camer. Comes to a stop and Indy rolls out"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 5
$ws.Range("C32").Value = 0.003472222222222222
$ws.Range("C32").NumberFormat = "h:mm"
$ws.Range("D32").Value = "MS"
$ws.Range("E32").Value = "Pan with Indy's movement"
$ws.Range("F32").Value = "Indy"
$ws.Range("G32").Value = "Indy finishes roll, catches breath and looks forward"

$ws.Range("A33").Value = 7
$ws.Range("B33").Value = 6
$ws.Range("C33").Value = 0.002777777777777778
$ws.Range("C33").NumberFormat = "h:mm"
$ws.Range("D33").Value = "MS, OTS"
$ws.Range("E33").Value = "Stationary"
$ws.Range("F33").Value = "Indy, Gopher"
$ws.Range("G33").Value = "Indy spots gopher. Gopher drops his snack and scurries into hole, while Indy stands up and walks off"

$ws.Range("A34").Value = 8
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 0.013194444444444444
$ws.Range("C34").NumberFormat = "h:mm"
$ws.Range("D34").Value = "MS/ES"
$ws.Range("E34").Value = "Cranes with Indy"
$ws.Range("F34").Value = "Indy, Nucluer explosion"
$ws.Range("G34").Value = "Indy walks up hill as a giant flash happens. He looks on as the mushroom cloud envelops the sky"

# -----------------------------------------------------------------
# 3. The "Stop about 3:08" note moves down to row 47.
# -----------------------------------------------------------------
$ws.Range("A47").Value = "Stop about 3:08"

# -----------------------------------------------------------------
# 4. Final selection, as left by the author.
# -----------------------------------------------------------------
$ws.Range("E37").Select() | Out-Null
